# tabel db.xlsx - perbaikan minor fungsi, testing bot ke line
# Adds a "timestamp/datetime" column and duplicates the two reference
# tables (tabel user / tabel utang) into a new block (columns K:R) while
# renaming a few fields to their Indonesian equivalents.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Rename a few existing field-name cells in "tabel utang" header rows
# ---------------------------------------------------------------------
$ws.Range("E12").Value = "id_borrower"
$ws.Range("G12").Value = "harga"
$ws.Range("H12").Value = "konfirmasi"

# ---------------------------------------------------------------------
# 2. New block titles first (this fixes the order new strings are
#    appended to the shared-string table: "Tabel User", "Tabel Utang",
#    then "datetime", "timestamp")
# ---------------------------------------------------------------------
$ws.Range("K2").Value = "Tabel User"
$ws.Range("K10").Value = "Tabel Utang"

# ---------------------------------------------------------------------
# 2b. New "Tabel User" block (mirrors B2:E8 into K2:N8, plus a new
#    "datetime"/"timestamp" column)
# ---------------------------------------------------------------------
$ws.Range("K3").Value = "int (auto incr)"
$ws.Range("L3").Value = "string (50)"
$ws.Range("M3").Value = "string (20)"
$ws.Range("N3").Value = "datetime"

$ws.Range("K4").Value = "id_user"
$ws.Range("L4").Value = "id_line"
$ws.Range("M4").Value = "username"
$ws.Range("N4").Value = "timestamp"

$ws.Range("K5").Value = 1
$ws.Range("K6").Value = 2
$ws.Range("K7").Value = 3
$ws.Range("K8").Value = 4

# ---------------------------------------------------------------------
# 3. New "Tabel Utang" block (mirrors C10:I16 into K10:R16, plus a new
#    "datetime"/"timestamp" column)
# ---------------------------------------------------------------------
$ws.Range("K11").Value = "int (auto incr)"
$ws.Range("L11").Value = "int"
$ws.Range("M11").Value = "int"
$ws.Range("N11").Value = "string (30)"
$ws.Range("O11").Value = "float"
$ws.Range("P11").Value = "bool"
$ws.Range("Q11").Value = "bool"
$ws.Range("R11").Value = "datetime"

$ws.Range("K12").Value = "nomor"
$ws.Range("L12").Value = "id_lender"
$ws.Range("M12").Value = "id_borrower"
$ws.Range("N12").Value = "komen"
$ws.Range("O12").Value = "harga"
$ws.Range("P12").Value = "konfirmasi"
$ws.Range("Q12").Value = "lunas"
$ws.Range("R12").Value = "timestamp"

$ws.Range("K13").Value = 1
$ws.Range("K14").Value = 2
$ws.Range("K15").Value = 3
$ws.Range("K16").Value = 4

# ---------------------------------------------------------------------
# 4. Formatting -- every box-bordered cell below copies its format from
#    an existing boxed cell (keeps reusing the workbook's existing thin
#    -border style instead of minting new ones). Source/target ranges
#    are always the same shape so Excel pastes a 1:1 cell mapping.
# ---------------------------------------------------------------------

# 4a. the "type" header row of "tabel utang" loses its box border
$ws.Range("C11:I11").Borders.LineStyle = -4142

# 4b. fill out row 16 of "tabel utang" with the same boxed-border look
#     as the rest of the table (copy format only, keep the cells' own
#     values/blankness)
$ws.Range("C13:I13").Copy()
$ws.Range("C16:I16").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 4c. new "Tabel User" block formatting (box border, matching C4:E8)
$ws.Range("C4:E4").Copy()
$ws.Range("K4:M4").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C4").Copy()
$ws.Range("N4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("C5:E7").Copy()
$ws.Range("K5:M7").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C5:C7").Copy()
$ws.Range("N5:N7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("C8:E8").Copy()
$ws.Range("K8:M8").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C4").Copy()
$ws.Range("N8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 4d. new "Tabel Utang" block formatting
$ws.Range("C11:I11").Copy()
$ws.Range("K11:Q11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("C12:I12").Copy()
$ws.Range("K12:Q12").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C12").Copy()
$ws.Range("R12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("C13:I13").Copy()
$ws.Range("K13:Q13").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C13").Copy()
$ws.Range("R13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("C14:I14").Copy()
$ws.Range("K14:Q14").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C14").Copy()
$ws.Range("R14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("C15:I15").Copy()
$ws.Range("K15:Q15").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C15").Copy()
$ws.Range("R15").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("C16:I16").Copy()
$ws.Range("K16:Q16").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C16").Copy()
$ws.Range("R16").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 5. Column widths / view
# ---------------------------------------------------------------------
$ws.Columns("E").AutoFit()
$ws.Columns("K").AutoFit()
$ws.Columns("M").AutoFit()
$ws.Columns("N").AutoFit()
$ws.Columns("R").AutoFit()

$ws.Range("M14").Select()
$excel.ActiveWindow.Zoom = 160
